# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the
# a2f32072-c7f3-4270-a7d9-182813d699b8.md row (row 5) in both the
# zh-cn and de-de localization-status sheets, reflecting a fresh handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-01 03:22:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-01 03:22:26"
